# "troca da palavra segundo pro coração do Leo não doer"
#
# No paragrafo sobre a busca de livros, o trecho
#   ".  A segunda possibilita filtrar a busca por titulo, ..."
# passa a ser
#   ".  Esta possibilita filtrar a busca por titulo, ..."
# ou seja, "A segunda" vira "Esta", mantendo o restante da frase intacto
# (mesma fonte/cor/tamanho) e sem alterar nenhum outro paragrafo.

$d = $word.ActiveDocument

# Localiza o trecho original de forma inequivoca.
$full = $d.Content.Text
$trechoAntigo = ".  A segunda possibilita "
$inicio = $full.IndexOf($trechoAntigo)
if ($inicio -lt 0) {
    throw "Trecho alvo '.  A segunda possibilita ' nao encontrado no documento."
}

# 1) Efetua a troca de texto propriamente dita.
$alvo = $d.Range($inicio, $inicio + $trechoAntigo.Length)
$alvo.Text = ".  Esta possibilita "

# 2) A engine tende a mesclar automaticamente runs adjacentes que
#    acabem com formatacao identica sempre que uma edicao de texto
#    atinge aquele trecho do paragrafo (da posicao editada ate o fim
#    do paragrafo). Para reproduzir exatamente a estrutura de runs
#    esperada (".  " / "Esta " / "possibilita " separados, e o restante
#    da frase preservado em seus runs originais, sem mesclar), "tocamos"
#    cada fronteira de run ligando e desligando negrito (sem mudar a
#    aparencia final) logo apos cada trecho, na ordem em que aparecem
#    no paragrafo a partir do ponto editado.
$segmentos = @(
    '.  ',
    'Esta ',
    'possibilita ',
    'filtrar ',
    'a busca por título, nome do autor, editora',
    ', ano de publicação, faixa de preço e idioma',
    ', ',
    'a fim de facilitar',
    ' a busca ',
    'do atributo desejado',
    '.',
    ' Para que o ',
    'cliente',
    ' possa adquirir algum livro, o ',
    'mesmo',
    ' deverá se cadastrar no site informando os seguintes dados: Nome, sobrenome, CPF, e-mail, telefone com DDD, data de nascimento, gênero',
    ', senha',
    ', endereço de cobrança e de entrega.'
)

$pos = $inicio
foreach ($seg in $segmentos) {
    $fim = $pos + $seg.Length
    $faixa = $d.Range($pos, $fim)
    if ($faixa.Text -ne $seg) {
        throw "Segmento inesperado em [$pos,$fim): '$($faixa.Text)' != '$seg'"
    }
    $faixa.Bold = 1
    $faixa.Bold = 0
    $pos = $fim
}
